# REVER_DailyTracker - "Add files via upload" update
# Adds new daily-tracker rows to several people's sheets, clears a stale
# entry on Rahman's sheet (re-aligning its row layout with the other
# sheets), and updates the active sheet/selection bookmarks.

$wb = $excel.ActiveWorkbook

$Balraj   = $wb.Worksheets.Item("Balraj")
$Vijay    = $wb.Worksheets.Item("Vijay")
$Bharathi = $wb.Worksheets.Item("Bharathi")
$Mamatha  = $wb.Worksheets.Item("Mamatha")
$Monisha  = $wb.Worksheets.Item("Monisha")
$Sabeena  = $wb.Worksheets.Item("Sabeena")
$Prabu    = $wb.Worksheets.Item("Prabu")
$Ram      = $wb.Worksheets.Item("Ram")
$Rahman   = $wb.Worksheets.Item("Rahman")
$Mathes   = $wb.Worksheets.Item("Mathes")

# ---------------------------------------------------------------------
# Balraj: new row 2 - RPA Dotnet / [Backend] CRM task, 30% WIP
# ---------------------------------------------------------------------
$Balraj.Range("A2").Value = 1
$Balraj.Range("B2").Value = 43950
$Balraj.Range("C2").Value = "RPA Dotnet"
$Balraj.Range("D2").Value = "[Backend] CRM"
$Balraj.Range("E2").Value = 0.3
$Balraj.Range("F2").Value = "WIP"
$Balraj.Range("G2").Value = "Following the discussion, we have to decided to have few more tables for the CRM which is in Progress, and also decided to  remove four columns at CusteorDetails and also at CustomerRequestData tables based on your suggestion."

# ---------------------------------------------------------------------
# Vijay: new row 2 - MujiStore authorization video task, 50% WIP
# ---------------------------------------------------------------------
$Vijay.Range("A2").Value = 1
$Vijay.Range("B2").Value = 43951
$Vijay.Range("C2").Value = "MujiStore"
$Vijay.Range("D2").Value = "MujiStore - Authorization for Stores video "
$Vijay.Range("E2").Value = 0.5
$Vijay.Range("F2").Value = "WIP"

# ---------------------------------------------------------------------
# Monisha: new row 2 - DotnetAPP Loginpage task, 100% Completed
# ---------------------------------------------------------------------
$Monisha.Range("A2").Value = 1
$Monisha.Range("B2").Value = 43950
$Monisha.Range("C2").Value = "DotnetAPP"
$Monisha.Range("D2").Value = "Loginpage (Frontend)"
$Monisha.Range("E2").Value = 1
$Monisha.Range("F2").Value = "Completed"

# ---------------------------------------------------------------------
# Sabeena: new rows 2 & 3 - RPA-SALE Video/PPT management tasks, WIP
# Their status cells pick up the "WIP" legend look (fill, no border),
# so copy that formatting from the legend row before writing values.
# ---------------------------------------------------------------------
$Vijay.Range("A2").Copy()
$Sabeena.Range("A2").PasteSpecial(-4122)
$Vijay.Range("B2").Copy()
$Sabeena.Range("B2").PasteSpecial(-4122)
$Vijay.Range("A2").Copy()
$Sabeena.Range("C2").PasteSpecial(-4122)
$Vijay.Range("A2").Copy()
$Sabeena.Range("D2").PasteSpecial(-4122)
$Vijay.Range("E2").Copy()
$Sabeena.Range("E2").PasteSpecial(-4122)
$Sabeena.Range("B23").Copy()
$Sabeena.Range("F2").PasteSpecial(-4122)
$Vijay.Range("A2").Copy()
$Sabeena.Range("G2").PasteSpecial(-4122)

$Vijay.Range("A2").Copy()
$Sabeena.Range("A3").PasteSpecial(-4122)
$Vijay.Range("B2").Copy()
$Sabeena.Range("B3").PasteSpecial(-4122)
$Vijay.Range("A2").Copy()
$Sabeena.Range("C3").PasteSpecial(-4122)
$Vijay.Range("A2").Copy()
$Sabeena.Range("D3").PasteSpecial(-4122)
$Vijay.Range("E2").Copy()
$Sabeena.Range("E3").PasteSpecial(-4122)
$Sabeena.Range("B23").Copy()
$Sabeena.Range("F3").PasteSpecial(-4122)
$Vijay.Range("A2").Copy()
$Sabeena.Range("G3").PasteSpecial(-4122)

$Sabeena.Range("A2").Value = 1
$Sabeena.Range("B2").Value = 43955
$Sabeena.Range("C2").Value = "RPA-SALE"
$Sabeena.Range("D2").Value = "Video Management - Create,view,update,delete"
$Sabeena.Range("E2").Value = 0.95
$Sabeena.Range("F2").Value = "WIP"

$Sabeena.Range("A3").Value = 2
$Sabeena.Range("B3").Value = 43956
$Sabeena.Range("C3").Value = "RPA-SALE"
$Sabeena.Range("D3").Value = "PPT Management - Create,view,update,delete"
$Sabeena.Range("E3").Value = 0.9
$Sabeena.Range("F3").Value = "WIP"

# ---------------------------------------------------------------------
# Rahman: clear the previously-entered "Excel macro program" row so it
# matches the other trackers' blank template row, and restore the row
# that had gone missing from the legend block (rows 19-26 realign with
# every other sheet in the workbook).
# ---------------------------------------------------------------------
$Rahman.Range("A2:G2").ClearContents()
$Rahman.Range("A2").Copy()
$Rahman.Range("D2").PasteSpecial(-4122)

$Rahman.Rows("19:19").Insert()

# ---------------------------------------------------------------------
# Selection / active-sheet bookmarks to mirror the saved workbook state
# ---------------------------------------------------------------------
$Balraj.Activate()
$Balraj.Range("E2").Select()

$Vijay.Activate()
$Vijay.Range("B2").Select()

$Monisha.Activate()
$Monisha.Range("A2:G2").Select()

$Rahman.Activate()
$Rahman.Range("F12").Select()

$Sabeena.Activate()
$Sabeena.Range("A3").Select()
